$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text cells (names, links, percentages, and multi-dot "prices") ---
$ws.Range("D2").Value = '26.307.74'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '1.678.56'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("E6").Value = '  +2.90%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("E8").Value = '  +2.17%  '
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("E10").Value = '  +1.42%  '
$ws.Range("E11").Value = '  +1.35%  '
$ws.Range("D12").Value = '1.687.06'
$ws.Range("E12").Value = '  +1.31%  '
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("E16").Value = '  +0.58%  '
$ws.Range("D17").Value = '26.333.27'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("E25").Value = '  +2.59%  '
$ws.Range("E26").Value = '  +3.49%  '
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("E28").Value = '  +0.83%  '
$ws.Range("E29").Value = '  +5.00%  '
$ws.Range("E30").Value = '  +0.85%  '
$ws.Range("E31").Value = '  +2.08%  '
$ws.Range("E32").Value = '  +2.29%  '
$ws.Range("E33").Value = '  +2.18%  '
$ws.Range("E34").Value = '  +1.07%  '
$ws.Range("E35").Value = '  +2.41%  '
$ws.Range("E36").Value = '  +1.96%  '
$ws.Range("E37").Value = '  +3.45%  '
$ws.Range("E38").Value = '  +3.87%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.108.47'
$ws.Range("E39").Value = '  +2.73%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E40").Value = '  +0.77%  '
$ws.Range("E41").Value = '  +2.07%  '
$ws.Range("E43").Value = '  -0.55%  '
$ws.Range("D44").Value = '1.830.50'
$ws.Range("E44").Value = '  +0.91%  '
$ws.Range("E45").Value = '  -5.49%  '
$ws.Range("E46").Value = '  +1.15%  '
$ws.Range("E47").Value = '  +1.31%  '
$ws.Range("E48").Value = '  -0.45%  '
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("E51").Value = '  +2.17%  '

# --- Numeric-looking price text that must remain stored as text ---
$numericTextCells = @{
    'D4' = '1.008'
    'D5' = '218.23'
    'D6' = '0.5262'
    'D8' = '0.2693'
    'D10' = '21.90'
    'D11' = '0.07523'
    'D13' = '4.523'
    'D14' = '0.5802'
    'D15' = '0.000008497'
    'D18' = '4.923'
    'D21' = '189.83'
    'D22' = '6.204'
    'D23' = '1.009'
    'D24' = '145.04'
    'D25' = '7.825'
    'D26' = '0.1249'
    'D27' = '15.78'
    'D28' = '0.06478'
    'D32' = '3.592'
    'D35' = '0.6228'
    'D36' = '2.408'
    'D38' = '6.413'
    'D40' = '0.01621'
    'D41' = '0.8765'
    'D43' = '100.47'
    'D45' = '0.00000000108'
    'D46' = '56.87'
    'D47' = '8.186'
    'D48' = '1.004'
    'D49' = '0.05270'
    'D51' = '6.075'
}
foreach ($ref in $numericTextCells.Keys) {
    $range = $ws.Range($ref)
    $range.NumberFormat = "@"
    $range.Value = $numericTextCells[$ref]
    $range.ClearFormats()
}
